$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Hour/Date/Epoch updated ---
# B2 "15:52:48" -> "21:05:41" (keep as text, not a Time value)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "21:05:41"
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats - restore original look (General/text)

# C2 "05-10-22" -> "06-10-22" (keep as text, not an auto-parsed date)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "06-10-22"
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# D2 epoch number
$ws.Range("D2").Value = 1665079541.348835

# --- Row 3: Hour/Date/Epoch updated ---
# B3 "15:52:53" -> "21:05:46"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "21:05:46"
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# C3 "05-10-22" -> "06-10-22"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "06-10-22"
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# D3 epoch number
$ws.Range("D3").Value = 1665079546.478611
